$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.668.23'
$ws.Range('E2').Value = '  -3.64%  '
$ws.Range('D3').Value = '2.627.61'
$ws.Range('E3').Value = '  -3.39%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''524.49'
$ws.Range('E5').Value = '  -1.01%  '
$ws.Range('D6').Value = '''143.35'
$ws.Range('E6').Value = '  -2.92%  '
$ws.Range('D7').Value = '''0.998'
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').Value = '''0.570'
$ws.Range('E8').Value = '  -1.63%  '
$ws.Range('E9').Value = '  -7.66%  '
$ws.Range('E10').Value = '  -2.52%  '
$ws.Range('D11').Value = '''0.336'
$ws.Range('E11').Value = '  -1.71%  '
$ws.Range('E12').Value = '  +0.91%  '
$ws.Range('D13').Value = '3.089.70'
$ws.Range('E13').Value = '  -3.28%  '
$ws.Range('D14').Value = '58.614.56'
$ws.Range('E14').Value = '  -3.69%  '
$ws.Range('D15').Value = '''21.08'
$ws.Range('E15').Value = '  -2.05%  '
$ws.Range('E16').Value = '  -1.70%  '
$ws.Range('D17').Value = '2.595.60'
$ws.Range('E17').Value = '  -6.17%  '
$ws.Range('D18').Value = '''338.45'
$ws.Range('E18').Value = '  -2.09%  '
$ws.Range('D19').Value = '''4.42'
$ws.Range('E19').Value = '  -2.10%  '
$ws.Range('D20').Value = '''10.45'
$ws.Range('E20').Value = '  -1.35%  '
$ws.Range('D21').Value = '''6.30'
$ws.Range('E21').Value = '  -2.34%  '
$ws.Range('D22').Value = '''1.00'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = '''65.16'
$ws.Range('E23').Value = '  +2.79%  '
$ws.Range('D24').Value = '''0.416'
$ws.Range('E24').Value = '  -0.77%  '
$ws.Range('D25').Value = '''0.167'
$ws.Range('D26').Value = '''0.997'
$ws.Range('E26').Value = '  +0.36%  '
$ws.Range('D27').Value = '''7.16'
$ws.Range('E27').Value = '  -2.41%  '
$ws.Range('D28').Value = '0.0₃0795'
$ws.Range('E28').Value = '  -3.85%  '
$ws.Range('D29').Value = '''6.51'
$ws.Range('E29').Value = '  -3.87%  '
$ws.Range('D30').Value = '''0.998'
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').Value = '''1.60'
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').Value = '''18.83'
$ws.Range('E32').Value = '  -1.32%  '
$ws.Range('D33').Value = '''150.35'
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').Value = '''4.13'
$ws.Range('E34').Value = '  -2.80%  '
$ws.Range('D35').Value = '''1.19'
$ws.Range('E35').Value = '  -3.41%  '
$ws.Range('D36').Value = '''0.902'
$ws.Range('E36').Value = '  -2.75%  '
$ws.Range('D37').Value = '''0.856'
$ws.Range('E37').Value = '  -5.40%  '
$ws.Range('D38').Value = '''36.49'
$ws.Range('E38').Value = '  -2.32%  '
$ws.Range('D39').Value = '''1.44'
$ws.Range('E39').Value = '  -6.44%  '
$ws.Range('D40').Value = '''3.63'
$ws.Range('E40').Value = '  -1.45%  '
$ws.Range('E41').Value = '  +0.33%  '
$ws.Range('D42').Value = '''0.604'
$ws.Range('E42').Value = '  -3.50%  '
$ws.Range('D43').Value = '''0.0972'
$ws.Range('E43').Value = '  -1.70%  '
$ws.Range('D44').Value = '''270.20'
$ws.Range('E44').Value = '  -4.18%  '
$ws.Range('D45').Value = '''10.66'
$ws.Range('E45').Value = '  +1.21%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '''0.0534'
$ws.Range('E46').Value = '  -1.83%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''19.14'
$ws.Range('E47').Value = '  -5.61%  '
$ws.Range('D48').Value = '2.034.86'
$ws.Range('E48').Value = '  -3.73%  '
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('D50').Value = '''4.59'
$ws.Range('E50').Value = '  -6.52%  '
$ws.Range('D51').Value = '''18.40'
$ws.Range('E51').Value = '  -5.50%  '
